# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column, and populate their data rows (2-26) with the values
# introduced by this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold / bordered / centered header style already used by the
# other header cells (copy format from H1, which carries that style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows ----------------------------------------------------------
$iValues = @(4,9,1,1,1,1,1,6,2,1,1,1,6,7,4,5,7,7,8,1,8,5,4,1,1)
$jValues = @(6,9,6,4,4,5,4,6,6,4,5,3,7,7,6,6,7,8,8,3,9,7,6,3,2)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Range("I$row").Value = $iValues[$k]
    $ws.Range("J$row").Value = $jValues[$k]
}

Write-Output "I0/IF columns added"
